$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @('-2.009878804933557','0.2541523497309064','-0.1970018051958322','-0.3373658687098063','0.1888731421524112','0.0470441954131372','0.6900188168152011','0.8144315231447835','0.4168013171071663','0.3406326034063205')
  3  = @('-1.893418122068486','-2.344572276995224','-2.484936340509199','-1.958697329646981','-2.100526276386255','-1.457551654984191','-1.333138948654609','-1.730769154692226','-1.806937868393072','-1.53061581027015')
  4  = @('1.095457709528105','0.9550936460141312','1.481332656876349','1.339503710137075','1.982478331539139','2.106891037868721','1.709260831831104','1.633092118130258','1.90941417625318','1.738194353941456')
  5  = @('1.016267824428332','1.54250683529055','1.400677888551276','2.04365250995334','2.168065216282922','1.770435010245305','1.694266296544459','1.970588354667381','1.799368532355657','1.952515937864398')
  6  = @('-0.4989491018657047','-0.6407780486049788','0.002196572797085183','0.1266092791266676','-0.2710209269109496','-0.3471896406117954','-0.07086758248887381','-0.2420874048005978','-0.08893999929185659','-0.3663687737149753')
  7  = @('-0.05818351157133772','0.5847911098307261','0.7092038161603086','0.3115736101226914','0.2354048964218456','0.5117269545447671','0.3405071322330432','0.4936545377417844','0.2162257633186657','0.05323806995971928')
  8  = @('0.4891160690213684','0.6135287753509509','0.2158985693133336','0.1397298556124878','0.4160519137354094','0.2448320914236854','0.3979794969324266','0.1205507225093079','-0.04243697084963852','-0.1665195462441563')
  9  = @('0.8739478123496736','0.4763176063120564','0.4001488926112106','0.6764709507341322','0.5052511284224082','0.6583985339311494','0.3809697595080307','0.2179820661490843','0.09389949075456649','0.6100525277605273')
  10 = @('-0.3591373898464348','-0.4353061035472806','-0.158984045424359','-0.330203867736083','-0.1770564622273418','-0.4544852366504605','-0.617472930009407','-0.7415555054039247','-0.2254024683979639','-0.4353061035472806')
  11 = @('-0.3352267436446591','-0.0589046855217375','-0.2301245078334615','-0.07697710232472027','-0.354405876747839','-0.5173935701067854','-0.6414761455013032','-0.1253231084953424','-0.3352267436446591',$null)
  12 = @('0.2948818205579588','0.1236619982462347','0.276809403754976','-0.0006193706681427817','-0.1636070640270892','-0.287689639421607','0.2284633975843539','0.01855976243503714',$null,$null)
  13 = @('-0.0245153921862106','0.1286320133225306','-0.1487967611005881','-0.3117844544595345','-0.4358670298540523','0.08028600715190851','-0.1296176279974082',$null,$null,$null)
  14 = @('-0.02881397568162436','-0.3062427501047431','-0.4692304434636895','-0.5933130188582073','-0.07715998185224648','-0.2870636170015632',$null,$null,$null,$null)
  15 = @('0.1944167064213277','0.0314290130623813','-0.09265356233213651','0.4234994746738243','0.2135958395245076',$null,$null,$null,$null,$null)
  16 = @('-0.2489288674730878','-0.3730114428676057','0.1431415941383551','-0.06676204101096155',$null,$null,$null,$null,$null,$null)
  17 = @('-0.201036585022594','0.3151164519833668','0.1052128168340501',$null,$null,$null,$null,$null,$null,$null)
  18 = @('0.009253912237035311','-0.2006497229122814',$null,$null,$null,$null,$null,$null,$null,$null)
  19 = @('0.4116802297750048',$null,$null,$null,$null,$null,$null,$null,$null,$null)
  20 = @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
}

$cols = @('B','C','D','E','F','G','H','I','J','K')

foreach ($r in ($data.Keys | Sort-Object)) {
  $rowVals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $val = $rowVals[$i]
    $cellAddr = "$col$r"
    if ($null -eq $val) {
      $ws.Range($cellAddr).Value = $null
    } else {
      $ws.Range($cellAddr).Value = [double]$val
    }
  }
}
